$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7103813333333333
$ws.Range("H2").Value = 2.131144
$ws.Range("I2").Value = 0.7576743564291667
$ws.Range("J2").Value = 0.7576743564291667
$ws.Range("M2").Value = 7.955277333333332
$ws.Range("N2").Value = 23.865832
$ws.Range("O2").Value = 0.05015625076675284
$ws.Range("P2").Value = 0.05015625076675283
$ws.Range("Q2").Value = 5.651280519089776
$ws.Range("R2").Value = 50.86152467180799
$ws.Range("S2").Value = 0.03800210502059936
$ws.Range("T2").Value = 0.03800210502059935
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7103813333333333
$ws.Range("H3").Value = 2.131144
$ws.Range("I3").Value = 0.7576743564291667
$ws.Range("J3").Value = 0.7576743564291667
$ws.Range("M3").Value = 82.48060333333333
$ws.Range("O3").Value = 0.520021823355633
$ws.Range("P3").Value = 0.520021823355633
$ws.Range("Q3").Value = 58.59268097007111
$ws.Range("R3").Value = 527.33412873064
$ws.Range("S3").Value = 0.394007200340101
$ws.Range("T3").Value = 0.394007200340101
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7103813333333333
$ws.Range("H4").Value = 2.131144
$ws.Range("I4").Value = 0.7576743564291667
$ws.Range("J4").Value = 0.7576743564291667
$ws.Range("M4").Value = 66.90297433333333
$ws.Range("N4").Value = 200.708923
$ws.Range("O4").Value = 0.4218083439585467
$ws.Range("P4").Value = 0.4218083439585465
$ws.Range("Q4").Value = 47.52662411087911
$ws.Range("R4").Value = 427.739616997912
$ws.Range("S4").Value = 0.3195933655452444
$ws.Range("T4").Value = 0.3195933655452444
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7103813333333333
$ws.Range("H5").Value = 2.131144
$ws.Range("I5").Value = 0.7576743564291667
$ws.Range("J5").Value = 0.7576743564291667
$ws.Range("M5").Value = 1.271033333333333
$ws.Range("N5").Value = 3.8131
$ws.Range("O5").Value = 0.008013581919067616
$ws.Range("P5").Value = 0.008013581919067614
$ws.Range("Q5").Value = 0.9029183540444445
$ws.Range("R5").Value = 8.126265186400001
$ws.Range("S5").Value = 0.006071685523221963
$ws.Range("T5").Value = 0.006071685523221961
$ws.Range("I6").Value = 0.2423256435708333
$ws.Range("J6").Value = 0.2423256435708333
$ws.Range("M6").Value = 7.955277333333332
$ws.Range("N6").Value = 23.865832
$ws.Range("O6").Value = 0.05015625076675284
$ws.Range("P6").Value = 0.05015625076675283
$ws.Range("Q6").Value = 1.807439010133333
$ws.Range("R6").Value = 16.2669510912
$ws.Range("S6").Value = 0.01215414574615348
$ws.Range("T6").Value = 0.01215414574615348
$ws.Range("I7").Value = 0.2423256435708333
$ws.Range("J7").Value = 0.2423256435708333
$ws.Range("M7").Value = 82.48060333333333
$ws.Range("O7").Value = 0.520021823355633
$ws.Range("P7").Value = 0.520021823355633
$ws.Range("S7").Value = 0.126014623015532
$ws.Range("T7").Value = 0.126014623015532
$ws.Range("I8").Value = 0.2423256435708333
$ws.Range("J8").Value = 0.2423256435708333
$ws.Range("M8").Value = 66.90297433333333
$ws.Range("N8").Value = 200.708923
$ws.Range("O8").Value = 0.4218083439585467
$ws.Range("P8").Value = 0.4218083439585465
$ws.Range("Q8").Value = 15.20035576853333
$ws.Range("R8").Value = 136.8032019168
$ws.Range("S8").Value = 0.1022149784133023
$ws.Range("T8").Value = 0.1022149784133022
$ws.Range("I9").Value = 0.2423256435708333
$ws.Range("J9").Value = 0.2423256435708333
$ws.Range("M9").Value = 1.271033333333333
$ws.Range("N9").Value = 3.8131
$ws.Range("O9").Value = 0.008013581919067616
$ws.Range("P9").Value = 0.008013581919067614
$ws.Range("Q9").Value = 0.2887787733333333
$ws.Range("R9").Value = 2.59900896
$ws.Range("S9").Value = 0.001941896395845654
$ws.Range("T9").Value = 0.001941896395845653
Write-Host "Updated NATMI LR-pair TPM values for Fgf8-Fgfr1 sheet"
